# fixed calculator bugs, rounded trend
# Adds a new row (row 15) of calculator/trend data to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

$ws.Cells.Item($row, 1).Value = 42622.890451388892
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = 18
$ws.Cells.Item($row, 3).Value = 57
$ws.Cells.Item($row, 4).Value = 42
$ws.Cells.Item($row, 5).Value = 57
$ws.Cells.Item($row, 6).Value = 18
$ws.Cells.Item($row, 7).Value = 30835
$ws.Cells.Item($row, 8).Value = 18162
$ws.Cells.Item($row, 9).Value = 978
$ws.Cells.Item($row, 10).Value = 163
$ws.Cells.Item($row, 11).Value = 120
$ws.Cells.Item($row, 12).Value = 22
$ws.Cells.Item($row, 13).Value = 5
$ws.Cells.Item($row, 14).Value = "Named"
